# Auto-generated edit script for horarios-141-completo workbook update
$wb = $excel.ActiveWorkbook

# ---- Sheet "LP1912" (sheet1) ----
$ws1 = $wb.Worksheets.Item("LP1912")
$ws1.Range("A2").Value = "Última actualización: 30/12/2025 10:28:34"
$ws1.Range("A3").Value = "Total filas: 130"

$sheet1Rows = @(
  @("10:28:23", "10:35", "23_HERNANDEZ", 7, "LP1912", "30/12/2025"),
  @("10:28:23", "10:42", "17_ROMERO", 14, "LP1912", "30/12/2025"),
  @("10:28:23", "10:43", "14_ABASTO", 15, "LP1912", "30/12/2025"),
  @("10:28:23", "10:57", "27_EL RETIRO", 29, "LP1912", "30/12/2025"),
  @("10:28:23", "11:02", "215C_EL PATO", 34, "LP1912", "30/12/2025"),
  @("10:28:23", "11:05", "23_HERNANDEZ", 37, "LP1912", "30/12/2025"),
  @("10:28:23", "11:06", "16_P MOR-167 Y 521", 38, "LP1912", "30/12/2025"),
  @("10:28:23", "11:11", "10_OLMOS", 43, "LP1912", "30/12/2025"),
  @("10:28:23", "11:21", "26_HERNANDEZ", 53, "LP1912", "30/12/2025"),
  @("10:28:23", "11:35", "23_HERNANDEZ", 67, "LP1912", "30/12/2025"),
  @("10:28:23", "11:42", "17_ROMERO", 74, "LP1912", "30/12/2025"),
  @("10:28:23", "11:52", "15_ABASTO", 84, "LP1912", "30/12/2025"),
  @("10:28:23", "11:53", "10_OLMOS", 85, "LP1912", "30/12/2025"),
  @("10:28:23", "12:02", "15_ABASTO", 94, "LP1912", "30/12/2025"),
  @("10:28:23", "12:02", "84_COLONIA URQUIZA-ESC 49", 94, "LP1912", "30/12/2025"),
  @("10:28:23", "12:03", "15_ABASTO", 95, "LP1912", "30/12/2025")
)

$r = 116
foreach ($row in $sheet1Rows) {
  $ws1.Cells.Item($r, 2).Value = $row[0]
  $ws1.Cells.Item($r, 3).Value = $row[1]
  $ws1.Cells.Item($r, 4).Value = $row[2]
  $ws1.Cells.Item($r, 5).Value = $row[3]
  $ws1.Cells.Item($r, 6).Value = $row[4]
  $ws1.Cells.Item($r, 7).Value = $row[5]
  $r = $r + 1
}

# ---- Sheet "LP1912-215" (sheet2) ----
$ws2 = $wb.Worksheets.Item("LP1912-215")
$ws2.Range("A2").Value = "Última actualización: 30/12/2025 10:28:34"
$ws2.Range("A3").Value = "Total filas: 16"
$ws2.Cells.Item(17, 2).Value = "30/12/2025"
$ws2.Cells.Item(17, 3).Value = "10:28:23"
$ws2.Cells.Item(17, 4).Value = "11:02"
$ws2.Cells.Item(17, 5).Value = "215C_EL PATO"
$ws2.Cells.Item(17, 6).Value = 34
$ws2.Cells.Item(17, 7).Value = "LP1912"

# ---- Sheet "6203-6173" (sheet3) ----
$ws3 = $wb.Worksheets.Item("6203-6173")
$ws3.Range("A2").Value = "Última actualización: 30/12/2025 10:28:34"
$ws3.Range("A3").Value = "Total filas: 19"

$sheet3Rows = @(
  @("30/12/2025", "10:28:34", "10:54", "215A_LA PLATA", 26, "L6173"),
  @("30/12/2025", "10:28:29", "11:14", "215C_LA PLATA", 46, "L6203"),
  @("30/12/2025", "10:28:34", "12:04", "215A_LA PLATA", 96, "L6173")
)

$r = 18
foreach ($row in $sheet3Rows) {
  $ws3.Cells.Item($r, 2).Value = $row[0]
  $ws3.Cells.Item($r, 3).Value = $row[1]
  $ws3.Cells.Item($r, 4).Value = $row[2]
  $ws3.Cells.Item($r, 5).Value = $row[3]
  $ws3.Cells.Item($r, 6).Value = $row[4]
  $ws3.Cells.Item($r, 7).Value = $row[5]
  $r = $r + 1
}
